# Update the "Förändrad" (C) date column for every data row, and add the
# friendly display-name second argument to the HYPERLINK() formulas that
# exist in the first block of rows (S, T, V, W, X, Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $firstRow + $used.Rows.Count - 1

$hyperlinkCols = @("S","T","V","W","X","Y")

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column C: bump the "changed" date serial from 45184 to 45186 ---
    $cAddr = "C" + $r
    $cVal = $ws.Range($cAddr).Value2()
    if ($cVal -eq 45184) {
        $ws.Range($cAddr).Value = 45186
    }

    # --- Hyperlink columns: append the display text argument ---
    $name = $ws.Range("A" + $r).Value()

    foreach ($col in $hyperlinkCols) {
        $addr = $col + $r
        $cell = $ws.Range($addr)
        $formula = $cell.Formula
        if ($formula -and $formula.Length -gt 0) {
            if ($formula.IndexOf("HYPERLINK(") -ge 0 -and $formula.IndexOf(",") -lt 0) {
                $trimmed = $formula.Substring(0, $formula.Length - 1)
                $cell.Formula = $trimmed + ', "' + $name + '")'
            }
        }
    }
}

Write-Host ("Updated rows 2.." + $lastRow)
